$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f18ce98f0de44b00a484f4a3cbcb3e17dcdcd921/e2e/e87024ab-88ca-4cf1-9f22-2d0e73c2e238.md"
$mdName = "e87024ab-88ca-4cf1-9f22-2d0e73c2e238.md"
$handedBackText = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $handedBackText
$ov.Range("F2").Value = $handedBackText
$ov.Columns.Item(5).ColumnWidth = 29.17
$ov.Columns.Item(6).ColumnWidth = 29.17

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $handedBackText
$zh.Columns.Item(3).ColumnWidth = 29.17
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

$zh.Range("I2").Value = $mdName
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, "", "", $mdName)
$zh.Range("I2").Font.Underline = 2
$zh.Range("I2").Font.Color = 15570276

$zh.Range("J2").Value = "e87024ab-88ca-4cf1-9f22-2d0e73c2e238.c440faa752b56b8e7ad176962d17787bd59bc3fc.zh-cn.xlf"
$zh.Range("K2").Value = "2016-10-19 17:34:35"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $handedBackText
$de.Columns.Item(3).ColumnWidth = 29.17
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17

$de.Range("I2").Value = $mdName
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, "", "", $mdName)
$de.Range("I2").Font.Underline = 2
$de.Range("I2").Font.Color = 15570276

$de.Range("J2").Value = "e87024ab-88ca-4cf1-9f22-2d0e73c2e238.c440faa752b56b8e7ad176962d17787bd59bc3fc.de-de.xlf"
$de.Range("K2").Value = "2016-10-19 17:34:53"
